$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the leave type name in C2 from "Paid Leave" to "Casual Leave"
$ws.Range("C2").Value = "Casual Leave"

# Reflect the active selection on the last-edited cell
$ws.Range("C2").Select()
